# Auto-generated: apply the per-cell text updates from the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to stage values that LOOK like plain numbers (e.g. "570.91")
# so Excel does not auto-convert them to the Number type. We stage the text in
# a helper cell formatted as Text, copy it, and Paste-Special "Values only" into
# the real target -- that brings across the literal string without carrying the
# helper cell's Text number-format onto the target (which must stay unstyled, as
# in the source file).
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

$ws.Range("D2").Value2 = "69.104.34"
$ws.Range("E2").Value2 = "  -2.20%  "
$ws.Range("D3").Value2 = "2.507.25"
$ws.Range("E3").Value2 = "  -0.73%  "
$ws.Range("E4").Value2 = "  +0.08%  "
$scratch.Value2 = "570.91"
$scratch.Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E5").Value2 = "  -0.70%  "
$scratch.Value2 = "166.10"
$scratch.Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E6").Value2 = "  -2.31%  "
$ws.Range("E7").Value2 = "  -0.04%  "
$ws.Range("E8").Value2 = "  +1.38%  "
$ws.Range("D9").Value2 = "2.507.50"
$ws.Range("E9").Value2 = "  -0.62%  "
$ws.Range("E10").Value2 = "  -2.19%  "
$ws.Range("E11").Value2 = "  -0.60%  "
$scratch.Value2 = "0.354"
$scratch.Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E12").Value2 = "  +3.05%  "
$scratch.Value2 = "4.90"
$scratch.Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E13").Value2 = "  +2.10%  "
$ws.Range("E14").Value2 = "  -0.28%  "
$ws.Range("D15").Value2 = "69.100.09"
$ws.Range("E15").Value2 = "  -1.86%  "
$scratch.Value2 = "0.0000174"
$scratch.Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E16").Value2 = "  -3.15%  "
$scratch.Value2 = "24.78"
$scratch.Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E17").Value2 = "  -0.37%  "
$ws.Range("D18").Value2 = "2.510.02"
$ws.Range("E18").Value2 = "  -0.84%  "
$scratch.Value2 = "11.33"
$scratch.Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E19").Value2 = "  -1.84%  "
$scratch.Value2 = "7.61"
$scratch.Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E20").Value2 = "  +0.93%  "
$scratch.Value2 = "348.65"
$scratch.Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E21").Value2 = "  -2.21%  "
$scratch.Value2 = "3.92"
$scratch.Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E22").Value2 = "  -0.17%  "
$ws.Range("E23").Value2 = "  +0.23%  "
$ws.Range("E24").Value2 = "  +0.11%  "
$scratch.Value2 = "70.26"
$scratch.Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E25").Value2 = "  +1.50%  "
$scratch.Value2 = "3.97"
$scratch.Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E26").Value2 = "  -1.89%  "
$ws.Range("E27").Value2 = "  -2.93%  "
$ws.Range("D28").Value2 = "2.637.67"
$ws.Range("E28").Value2 = "  -0.41%  "
$scratch.Value2 = "0.999"
$scratch.Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E29").Value2 = "  +0.06%  "
$ws.Range("D30").Value2 = "0.0₃0889"
$ws.Range("E30").Value2 = "  -2.22%  "
$scratch.Value2 = "7.85"
$scratch.Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E31").Value2 = "  -0.02%  "
$scratch.Value2 = "460.63"
$scratch.Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E32").Value2 = "  -3.80%  "
$ws.Range("E33").Value2 = "  -1.89%  "
$ws.Range("E34").Value2 = "  -1.65%  "
$ws.Range("E35").Value2 = "  +0.20%  "
$ws.Range("E36").Value2 = "  +1.58%  "
$scratch.Value2 = "157.10"
$scratch.Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E37").Value2 = "  -0.24%  "
$scratch.Value2 = "18.99"
$scratch.Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$scratch.Value2 = "18.49"
$scratch.Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E39").Value2 = "  -0.54%  "
$scratch.Value2 = "4.74"
$scratch.Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E41").Value2 = "  +0.59%  "
$ws.Range("E42").Value2 = "  -0.35%  "
$ws.Range("E43").Value2 = "  -3.14%  "
$scratch.Value2 = "38.22"
$scratch.Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E44").Value2 = "  -0.17%  "
$ws.Range("E45").Value2 = "  -13.19%  "
$scratch.Value2 = "2.24"
$scratch.Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E46").Value2 = "  -6.55%  "
$scratch.Value2 = "141.27"
$scratch.Copy() | Out-Null
$ws.Range("D47").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E47").Value2 = "  -1.30%  "
$ws.Range("E48").Value2 = "  +0.32%  "
$ws.Range("E49").Value2 = "  -1.57%  "
$scratch.Value2 = "0.0729"
$scratch.Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E50").Value2 = "  -0.84%  "
$scratch.Value2 = "1.56"
$scratch.Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("E51").Value2 = "  -3.70%  "

# Tidy up: remove the scratch cell and clipboard marquee entirely so the
# sheet dimensions / used range stay exactly as before (A1:E51).
$scratch.Clear()
$excel.CutCopyMode = 0
